# Update NATMI ligand-receptor pair metrics with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.027114666666666
$ws.Range("H2").Value = 9.081344
$ws.Range("I2").Value = 0.207506525262911
$ws.Range("J2").Value = 0.207506525262911
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2758483333333333
$ws.Range("N2").Value = 0.827545
$ws.Range("Q2").Value = 0.8350245356088888
$ws.Range("R2").Value = 7.51522082048
$ws.Range("S2").Value = 0.207506525262911
$ws.Range("T2").Value = 0.207506525262911

# Row 3
$ws.Range("G3").Value = 2.518570333333333
$ws.Range("H3").Value = 7.555711000000001
$ws.Range("I3").Value = 0.1726461783080517
$ws.Range("J3").Value = 0.1726461783080516
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2758483333333333
$ws.Range("N3").Value = 0.827545
$ws.Range("Q3").Value = 0.6947434288327777
$ws.Range("R3").Value = 6.252690859495
$ws.Range("S3").Value = 0.1726461783080517
$ws.Range("T3").Value = 0.1726461783080516

# Row 4
$ws.Range("G4").Value = 4.235286666666666
$ws.Range("H4").Value = 12.70586
$ws.Range("I4").Value = 0.2903258437382188
$ws.Range("J4").Value = 0.2903258437382187
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2758483333333333
$ws.Range("N4").Value = 0.827545
$ws.Range("Q4").Value = 1.168296768188889
$ws.Range("R4").Value = 10.5146709137
$ws.Range("S4").Value = 0.2903258437382188
$ws.Range("T4").Value = 0.2903258437382187

# Row 5
$ws.Range("G5").Value = 1.937427333333333
$ws.Range("H5").Value = 5.812282
$ws.Range("I5").Value = 0.1328092451588843
$ws.Range("J5").Value = 0.1328092451588843
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2758483333333333
$ws.Range("N5").Value = 0.827545
$ws.Range("Q5").Value = 0.5344361008544444
$ws.Range("R5").Value = 4.809924907689999
$ws.Range("S5").Value = 0.1328092451588843
$ws.Range("T5").Value = 0.1328092451588843

# Row 6
$ws.Range("G6").Value = 2.869646666666667
$ws.Range("H6").Value = 8.60894
$ws.Range("I6").Value = 0.1967122075319342
$ws.Range("J6").Value = 0.1967122075319342
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2758483333333333
$ws.Range("N6").Value = 0.827545
$ws.Range("Q6").Value = 0.7915872502555554
$ws.Range("R6").Value = 7.1242852523
$ws.Range("S6").Value = 0.1967122075319342
$ws.Range("T6").Value = 0.1967122075319342
